# Add the header row for the vaccinated-students sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen column E (Vaccine_Dose) so the header isn't truncated.
$ws.Columns.Item(5).ColumnWidth = 13.3

# Leave the selection where the author's cursor ended up after typing
# the last header (one cell to the right, on F1).
[void]$ws.Range("F1").Select()
